$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.250.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.651.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3900'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.33%  '
$ws.Range("E8").Value = '  -3.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.004'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("E10").Value = '  -7.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '49.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08506'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.183'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001294'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.548'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.650.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.04'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06915'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.74%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '21.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.988'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.96%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("E23").Value = '  -4.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.242.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.368'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.770'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '158.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.558'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '143.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.385'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -13.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.434'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.830.71'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.973'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08157'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9979'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02943'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2728'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09339'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.490'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7688'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.517'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6928'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.107'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.38%  '
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08464'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.277'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '135.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.99%  '
